# B6-PowerPoint.pptx edit:
#  1. Re-colour the deck's theme from the "Integral" (Red Violet) colour
#     values to the standard "Office" colour values.
#  2. Re-style the three tables in the deck that used the old custom
#     table style with the new table style id.

function HexToRgbInt($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# New ("Office") theme colour scheme, in clrScheme order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = HexToRgbInt $officeColors[$i - 1]
}

# Swap every table still using the old custom table style for the new one.
$oldStyleId = "{58B64D4B-4D15-4EE1-810A-E03B1601406B}"
$newStyleId = "{D01A47BA-43DE-4A8D-8010-E5754F0CC432}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            if ($shape.Table.Style -eq $oldStyleId) {
                $shape.Table.ApplyStyle($newStyleId)
            }
        }
    }
}
